$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44335
$ws.Range("M2").Value = 80
$ws.Range("N2").Value = 14000
$ws.Range("O2").Value = 14000
$ws.Range("P2").Value = 14000
$ws.Range("R2").Value = "Región Metropolitana"
$ws.Range("S2").Value = 2000
$ws.Range("D3").Value = 44307
$ws.Range("L3").Value = "Primera"
$ws.Range("N3").Value = 14000
$ws.Range("O3").Value = 14000
$ws.Range("P3").Value = 14000
$ws.Range("S3").Value = 2000
$ws.Range("D4").Value = 44307
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 10000
$ws.Range("O4").Value = 10000
$ws.Range("P4").Value = 10000
$ws.Range("S4").Value = 1429
$ws.Range("D5").Value = 44314
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 20
$ws.Range("N5").Value = 13000
$ws.Range("O5").Value = 13000
$ws.Range("P5").Value = 13000
$ws.Range("S5").Value = 1857
$ws.Range("D6").Value = 44314
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 45
$ws.Range("N6").Value = 11000
$ws.Range("O6").Value = 11000
$ws.Range("P6").Value = 11000
$ws.Range("S6").Value = 1571
$ws.Range("D7").Value = 44342
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 12000
$ws.Range("O7").Value = 12000
$ws.Range("P7").Value = 12000
$ws.Range("S7").Value = 1714
$ws.Range("D8").Value = 44300
$ws.Range("M8").Value = 150
$ws.Range("O8").Value = 13000
$ws.Range("P8").Value = 12500
$ws.Range("R8").Value = "Provincia de Santiago"
$ws.Range("S8").Value = 1786
$ws.Range("D9").Value = 44315
$ws.Range("L9").Value = "Especial"
$ws.Range("M9").Value = 50
$ws.Range("N9").Value = 14000
$ws.Range("O9").Value = 14000
$ws.Range("P9").Value = 14000
$ws.Range("S9").Value = 2000
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 80
$ws.Range("N10").Value = 12000
$ws.Range("O10").Value = 13000
$ws.Range("P10").Value = 12500
$ws.Range("S10").Value = 1786
$ws.Range("L11").Value = "Segunda"
$ws.Range("N11").Value = 10000
$ws.Range("O11").Value = 11000
$ws.Range("P11").Value = 10500
$ws.Range("S11").Value = 1500
$ws.Range("D12").Value = 44312
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 50
$ws.Range("N12").Value = 13000
$ws.Range("O12").Value = 13000
$ws.Range("P12").Value = 13000
$ws.Range("S12").Value = 1857
$ws.Range("D13").Value = 44312
$ws.Range("L13").Value = "Segunda"
$ws.Range("M13").Value = 20
$ws.Range("D14").Value = 44349
$ws.Range("L14").Value = "Segunda"
$ws.Range("M14").Value = 70
$ws.Range("O14").Value = 12000
$ws.Range("P14").Value = 12000
$ws.Range("R14").Value = "Región Metropolitana"
$ws.Range("S14").Value = 1714
$ws.Range("D15").Value = 44316
$ws.Range("M15").Value = 40
$ws.Range("N15").Value = 13000
$ws.Range("O15").Value = 13000
$ws.Range("P15").Value = 13000
$ws.Range("S15").Value = 1857
$ws.Range("D16").Value = 44316
$ws.Range("N16").Value = 11000
$ws.Range("O16").Value = 11000
$ws.Range("P16").Value = 11000
$ws.Range("S16").Value = 1571
$ws.Range("D17").Value = 44306
$ws.Range("M17").Value = 50
$ws.Range("N17").Value = 12000
$ws.Range("O17").Value = 12000
$ws.Range("P17").Value = 12000
$ws.Range("S17").Value = 1714
$ws.Range("D18").Value = 44306
$ws.Range("M18").Value = 40
$ws.Range("N18").Value = 9000
$ws.Range("O18").Value = 9000
$ws.Range("P18").Value = 9000
$ws.Range("S18").Value = 1286
$ws.Range("D19").Value = 44322
$ws.Range("M19").Value = 100
$ws.Range("N19").Value = 11000
$ws.Range("O19").Value = 11000
$ws.Range("P19").Value = 11000
$ws.Range("S19").Value = 1571
$ws.Range("D20").Value = 44321
$ws.Range("M20").Value = 140
$ws.Range("N20").Value = 11000
$ws.Range("O20").Value = 12000
$ws.Range("P20").Value = 11500
$ws.Range("S20").Value = 1643
$ws.Range("D21").Value = 44321
$ws.Range("M21").Value = 80
$ws.Range("N21").Value = 8000
$ws.Range("O21").Value = 8000
$ws.Range("P21").Value = 8000
$ws.Range("S21").Value = 1143
$ws.Range("D22").Value = 44302
$ws.Range("L22").Value = "Primera"
$ws.Range("M22").Value = 340
$ws.Range("O22").Value = 13000
$ws.Range("P22").Value = 12500
$ws.Range("R22").Value = "Provincia de Santiago"
$ws.Range("S22").Value = 1786
